$d = $word.ActiveDocument

# --- Fix footnote #2 indentation: add a tab-stop leader + leading tab before
# the footnote reference mark, and a tab between the reference mark and the
# footnote text (instead of a literal leading space in the text run).
$footnoteXml = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:pPr>
<w:tabs>
<w:tab w:val="start" w:pos="0"/>
<w:tab w:val="start" w:pos="400"/>
</w:tabs>
<w:suppressAutoHyphens/>
</w:pPr>
<w:r><w:tab/></w:r>
<w:r><w:rPr><w:vertAlign w:val="superscript"/></w:rPr><w:footnoteRef/></w:r>
<w:r><w:tab/></w:r>
<w:r><w:t xml:space="preserve">This is the footnote content.</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@

$d.Footnotes.Item(1).Range.InsertXML($footnoteXml)

# --- Same fix for endnote #2 ---
$endnoteXml = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:pPr>
<w:tabs>
<w:tab w:val="start" w:pos="0"/>
<w:tab w:val="start" w:pos="400"/>
</w:tabs>
<w:suppressAutoHyphens/>
</w:pPr>
<w:r><w:tab/></w:r>
<w:r><w:rPr><w:vertAlign w:val="superscript"/></w:rPr><w:endnoteRef/></w:r>
<w:r><w:tab/></w:r>
<w:r><w:t xml:space="preserve">This is the endnote content.</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@

$d.Endnotes.Item(1).Range.InsertXML($endnoteXml)
